$wb = $excel.ActiveWorkbook

# --- Performance Metrics sheet: remove extraneous space after comma in ranges ---
$wsPerf = $wb.Worksheets.Item("Performance Metrics")
$wsPerf.Range("J2").Value = "1.24 [1.2,1.29]"
$wsPerf.Range("N2").Value = "0.581 [0.571,0.592]"

# --- Evaluation Sample Sets sheet: switch dict-style strings to the new REST API v1.4 format ---
$wsEval = $wb.Worksheets.Item("Evaluation Sample Sets")
$wsEval.Range("G2").Value = "mean:54.3;range:[50.1,58.4];unit:years"
$wsEval.Range("M2").Value = "mean:6.2;sd:1.7;unit:years"
